$d = $word.ActiveDocument

$pairs = @(
    @{ old = "2025-07-31 Thursday"; new = "2025-08-01 Friday" },
    @{ old = "641×7=4487"; new = "224×9=2016" },
    @{ old = "298×8=2384"; new = "954×5=4770" },
    @{ old = "400×9=3600"; new = "560×5=2800" },
    @{ old = "291×5=1455"; new = "681×9=6129" },
    @{ old = "733×9=6597"; new = "217×9=1953" },
    @{ old = "543×6=3258"; new = "510×8=4080" },
    @{ old = "255×3=765"; new = "430×8=3440" },
    @{ old = "913×5=4565"; new = "599×6=3594" },
    @{ old = "104×8=832"; new = "850×9=7650" },
    @{ old = "133×3=399"; new = "102×9=918" },
    @{ old = "266×4=1064"; new = "467×9=4203" },
    @{ old = "848×6=5088"; new = "432×7=3024" },
    @{ old = "461×7=3227"; new = "906×8=7248" },
    @{ old = "795×5=3975"; new = "889×8=7112" },
    @{ old = "576×7=4032"; new = "562×8=4496" },
    @{ old = "147×2=294"; new = "561×9=5049" },
    @{ old = "382×4=1528"; new = "928×5=4640" },
    @{ old = "493×3=1479"; new = "220×3=660" },
    @{ old = "833×9=7497"; new = "590×6=3540" },
    @{ old = "908×9=8172"; new = "385×6=2310" },
    @{ old = "228×4=912"; new = "730×6=4380" },
    @{ old = "113×3=339"; new = "368×8=2944" },
    @{ old = "682×4=2728"; new = "585×9=5265" },
    @{ old = "748×4=2992"; new = "389×5=1945" },
    @{ old = "558×6=3348"; new = "114×6=684" }
)

foreach ($pair in $pairs) {
    $range = $d.Content
    $range.Find.Execute($pair.old, $true, $false, $false, $false, $false, $true, 1, $false, $pair.new, 2)
}
